# Refresh the "cryptos" price/volume snapshot (GitHub Actions daily update).
# Updates Price (D) / Volume(1h) (E) text for each coin row, and for the two
# rows whose ranking swapped (ImmutableX<->Dai, Kaspa<->Stellar) also updates
# Coin (B) and Link (C). Values that look like plain numbers are written with
# a leading apostrophe so they stay text cells, matching the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.523.99'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '2.462.67'
$ws.Range('E3').Value = '  -0.70%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.66%  '
$ws.Range('D5').Value = '''313.42'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').Value = '''91.17'
$ws.Range('E6').Value = '  -2.16%  '
$ws.Range('D7').Value = '''0.548'
$ws.Range('E7').Value = '  +0.31%  '
$ws.Range('E8').Value = '  -0.53%  '
$ws.Range('E9').Value = '  +3.51%  '
$ws.Range('D10').Value = '''32.39'
$ws.Range('E10').Value = '  -2.62%  '
$ws.Range('D11').Value = '''0.0790'
$ws.Range('E11').Value = '  +1.31%  '
$ws.Range('E12').Value = '  +0.65%  '
$ws.Range('D13').Value = '2.845.22'
$ws.Range('E13').Value = '  -0.83%  '
$ws.Range('D14').Value = '''6.82'
$ws.Range('E14').Value = '  -0.84%  '
$ws.Range('D15').Value = '''15.82'
$ws.Range('E15').Value = '  +2.85%  '
$ws.Range('D16').Value = '2.457.99'
$ws.Range('E16').Value = '  -1.02%  '
$ws.Range('D17').Value = '''0.773'
$ws.Range('E17').Value = '  -1.72%  '
$ws.Range('D18').Value = '41.489.58'
$ws.Range('E18').Value = '  +0.26%  '
$ws.Range('D19').Value = '''6.49'
$ws.Range('E19').Value = '  +3.02%  '
$ws.Range('D20').Value = '0.0₃0939'
$ws.Range('E20').Value = '  +1.42%  '
$ws.Range('D21').Value = '''70.73'
$ws.Range('E21').Value = '  +1.00%  '
$ws.Range('D22').Value = '''11.03'
$ws.Range('E22').Value = '  -0.82%  '
$ws.Range('D23').Value = '''237.27'
$ws.Range('E23').Value = '  +1.03%  '
$ws.Range('E24').Value = '  -1.44%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = '''1.00'
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').Value = '''1.90'
$ws.Range('E26').Value = '  +0.98%  '
$ws.Range('D27').Value = '''24.50'
$ws.Range('E27').Value = '  +1.73%  '
$ws.Range('E28').Value = '  -0.48%  '
$ws.Range('E29').Value = '  -1.24%  '
$ws.Range('D30').Value = '''35.20'
$ws.Range('E30').Value = '  -3.63%  '
$ws.Range('D31').Value = '''155.55'
$ws.Range('E31').Value = '  +1.77%  '
$ws.Range('D32').Value = '''5.41'
$ws.Range('E32').Value = '  -0.92%  '
$ws.Range('E33').Value = '  +0.64%  '
$ws.Range('D34').Value = '''0.0755'
$ws.Range('E34').Value = '  +0.60%  '
$ws.Range('D35').Value = '''17.09'
$ws.Range('E35').Value = '  -3.71%  '
$ws.Range('E36').Value = '  -6.32%  '
$ws.Range('D37').Value = '''2.86'
$ws.Range('E37').Value = '  -5.48%  '
$ws.Range('B38').Value = 'Stellar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D38').Value = '''0.114'
$ws.Range('E38').Value = '  +0.88%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = '''0.103'
$ws.Range('E39').Value = '  +2.77%  '
$ws.Range('D40').Value = '''1.77'
$ws.Range('E40').Value = '  -4.36%  '
$ws.Range('E41').Value = '  -2.34%  '
$ws.Range('E42').Value = '  -0.73%  '
$ws.Range('D43').Value = '1.940.29'
$ws.Range('E43').Value = '  -1.80%  '
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('D45').Value = '''18.59'
$ws.Range('E45').Value = '  -5.34%  '
$ws.Range('E46').Value = '  -3.01%  '
$ws.Range('D47').Value = '''9.00'
$ws.Range('E47').Value = '  +2.41%  '
$ws.Range('D48').Value = '2.707.53'
$ws.Range('E48').Value = '  -0.97%  '
$ws.Range('D49').Value = '''96.58'
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('E50').Value = '  -2.67%  '
$ws.Range('D51').Value = '''52.05'
$ws.Range('E51').Value = '  +3.36%  '
